$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LUAD-bic")
$ws.Range("A2").Value = "'16"
$ws.Range("B2").Value = "Mutation EGFR"
$ws.Range("C2").Value = "Mutation PIK3CA"
$ws.Range("E2").Value = [double]"1.0"
$ws.Range("F2").Value = [double]"4.23751999471512e-10"
$ws.Range("G2").Value = [double]"2.46328156654618e-10"
$ws.Range("H2").Value = [double]"0.0"
$ws.Range("K2").Value = [double]"0.13125"
$ws.Range("L2").Value = [double]"0.03547788826234667"
$ws.Range("M2").Value = [double]"0.15"
$ws.Range("N2").Value = [double]"0.0437003686737563"
$ws.Range("A3").Value = "'17"
$ws.Range("B3").Value = "Mutation NF1"
$ws.Range("C3").Value = "Mutation RB1"
$ws.Range("E3").Value = [double]"2.0"
$ws.Range("F3").Value = [double]"0.000498638730127306"
$ws.Range("G3").Value = [double]"1.56420124039081e-09"
$ws.Range("H3").Value = [double]"0.025"
$ws.Range("K3").Value = [double]"0.09375"
$ws.Range("L3").Value = [double]"0.03294039229342062"
$ws.Range("M3").Value = [double]"0.14375"
$ws.Range("N3").Value = [double]"0.0514545376467853"
$ws.Range("K4").Value = [double]"0.125"
$ws.Range("L4").Value = [double]"0.0"
$ws.Range("A5").Value = "'7"
$ws.Range("B5").Value = "Amplification RIT1"
$ws.Range("C5").Value = "Mutation AKT1"
$ws.Range("D5").Value = [double]"2.0"
$ws.Range("F5").Value = [double]"0.000511067883527512"
$ws.Range("G5").Value = [double]"4.21328460987917e-08"
$ws.Range("K5").Value = [double]"0.125"
$ws.Range("L5").Value = [double]"0.0"
$ws.Range("N5").Value = [double]"0.0197642353760524"
$ws.Range("A6").Value = "'9"
$ws.Range("B6").Value = "Mutation ALK"
$ws.Range("C6").Value = "Mutation ARID1A"
$ws.Range("D6").Value = [double]"3.0"
$ws.Range("F6").Value = [double]"2.497294669101e-09"
$ws.Range("G6").Value = [double]"3.31352455569443e-09"
$ws.Range("I6").Value = [double]"40.0"
$ws.Range("K6").Value = [double]"0.14375"
$ws.Range("L6").Value = [double]"0.05145453764678529"
$ws.Range("M6").Value = [double]"0.15"
$ws.Range("N6").Value = [double]"0.0322748612183951"
$ws.Range("I7").Value = [double]"30.0"
$ws.Range("K8").Value = [double]"0.125"
$ws.Range("L8").Value = [double]"0.0"
$ws.Range("K25").Value = [double]"0.575"
$ws.Range("L25").Value = [double]"0.09682458365518543"

$ws = $wb.Worksheets.Item("LUAD-aic")
$ws.Range("A2").Value = "'17"
$ws.Range("B2").Value = "Mutation EGFR"
$ws.Range("C2").Value = "Mutation PIK3CA"
$ws.Range("E2").Value = [double]"1.0"
$ws.Range("F2").Value = [double]"4.23751999471512e-10"
$ws.Range("G2").Value = [double]"2.46328156654618e-10"
$ws.Range("H2").Value = [double]"0.0"
$ws.Range("K2").Value = [double]"0.13125"
$ws.Range("L2").Value = [double]"0.03547788826234667"
$ws.Range("M2").Value = [double]"0.1125"
$ws.Range("N2").Value = [double]"0.0395284707521047"
$ws.Range("A3").Value = "'18"
$ws.Range("B3").Value = "Mutation NF1"
$ws.Range("C3").Value = "Mutation RB1"
$ws.Range("E3").Value = [double]"2.0"
$ws.Range("F3").Value = [double]"0.000498638730127306"
$ws.Range("G3").Value = [double]"1.56420124039081e-09"
$ws.Range("H3").Value = [double]"0.025"
$ws.Range("K3").Value = [double]"0.09375"
$ws.Range("L3").Value = [double]"0.03294039229342062"
$ws.Range("M3").Value = [double]"0.15"
$ws.Range("N3").Value = [double]"0.0322748612183951"
$ws.Range("K4").Value = [double]"0.125"
$ws.Range("L4").Value = [double]"0.0"
$ws.Range("A5").Value = "'8"
$ws.Range("B5").Value = "Amplification RIT1"
$ws.Range("C5").Value = "Mutation AKT1"
$ws.Range("D5").Value = [double]"2.0"
$ws.Range("F5").Value = [double]"0.000511067883527512"
$ws.Range("G5").Value = [double]"4.21328460987917e-08"
$ws.Range("K5").Value = [double]"0.125"
$ws.Range("L5").Value = [double]"0.0"
$ws.Range("N5").Value = [double]"0.0"
$ws.Range("A6").Value = "'10"
$ws.Range("B6").Value = "Mutation ALK"
$ws.Range("C6").Value = "Mutation ARID1A"
$ws.Range("D6").Value = [double]"3.0"
$ws.Range("F6").Value = [double]"2.497294669101e-09"
$ws.Range("G6").Value = [double]"3.31352455569443e-09"
$ws.Range("I6").Value = [double]"40.0"
$ws.Range("K6").Value = [double]"0.14375"
$ws.Range("L6").Value = [double]"0.05145453764678529"
$ws.Range("M6").Value = [double]"0.1375"
$ws.Range("N6").Value = [double]"0.0493006648591635"
$ws.Range("I7").Value = [double]"30.0"
$ws.Range("K8").Value = [double]"0.125"
$ws.Range("L8").Value = [double]"0.0"
$ws.Range("M8").Value = [double]"0.125"
$ws.Range("N8").Value = [double]"0.0"
$ws.Range("M11").Value = [double]"0.18125"
$ws.Range("N11").Value = [double]"0.0197642353760524"
$ws.Range("K26").Value = [double]"0.575"
$ws.Range("L26").Value = [double]"0.09682458365518543"
